$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store plain numeric-looking and
# percentage-looking text as literal strings (e.g. "276.13", "-1.02%").
# Format each touched cell as Text first so Excel keeps the new value as a
# literal string instead of auto-converting it to a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.02%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.24%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.874"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.57%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.19%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.912"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.55%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.321"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.34%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.303"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "37.05%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8751"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.25%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1565"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.39%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05010"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.62%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07467"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.88%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02977"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-5.02%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09058"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.12%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001573"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.94%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006307"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.59%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006056"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.85%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.26%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.32%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.15%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1333"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.62%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.915"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.25%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04369"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.22%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001174"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.51%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004211"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.06%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.00%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-4.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04098"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.20%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007001"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.67%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1173"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.71%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.30%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01119"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.43%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005302"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.47%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.02001"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-11.12%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.490"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.35%"
